# Daily update at 8 AM UTC
# Appends the next day's win counts to the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

$ws.Cells.Item($row, 1).Value = 45983
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = 75
$ws.Cells.Item($row, 3).Value = 86
$ws.Cells.Item($row, 4).Value = 81
